# Update Name of Algo
# Applies updated RandomForest imputation results to column D and E cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.9387
$ws.Range("D4").Value = -7.824500000000001
$ws.Range("E6").Value = 12.622
$ws.Range("D7").Value = -8.433399999999999
$ws.Range("E7").Value = 11.8324
$ws.Range("D8").Value = -8.4901
$ws.Range("E8").Value = 12.34920000000001
$ws.Range("B11").Value = 5.027299999999996
$ws.Range("B12").Value = 4.894699999999996
$ws.Range("D12").Value = -8.006399999999999
$ws.Range("D14").Value = -8.570599999999999
$ws.Range("B15").Value = 4.850600000000002
$ws.Range("E19").Value = 13.0103
$ws.Range("E21").Value = 12.7154
$ws.Range("D22").Value = -8.003199999999994
$ws.Range("E24").Value = 12.7972
$ws.Range("E25").Value = 13.2542
